$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Dingo"
$ws.Range("B8").Value = 12

$ws.Range("A9").Value = "Gerbbbbb"
$ws.Range("B9").Value = 12
